$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.565.72'
$ws.Range("E2").Value = '  +7.15%  '
$ws.Range("D3").Value = '1.725.83'
$ws.Range("E3").Value = '  +3.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '332.20'
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9969'
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3710'
$ws.Range("E7").Value = '  +1.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.51'
$ws.Range("E8").Value = '  +2.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3362'
$ws.Range("E9").Value = '  +3.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.185'
$ws.Range("E10").Value = '  +4.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07401'
$ws.Range("E11").Value = '  +4.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9962'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.387'
$ws.Range("E13").Value = '  +5.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.11'
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.061'
$ws.Range("E15").Value = '  +7.13%  '
$ws.Range("D16").Value = '1.724.21'
$ws.Range("E16").Value = '  +3.66%  '
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06630'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.89'
$ws.Range("E19").Value = '  +4.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9981'
$ws.Range("E20").Value = '  -0.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.56'
$ws.Range("E21").Value = '  +5.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.135'
$ws.Range("E22").Value = '  +3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.76'
$ws.Range("E23").Value = '  +1.95%  '
$ws.Range("D24").Value = '26.569.94'
$ws.Range("E24").Value = '  +7.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.436'
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.424'
$ws.Range("E26").Value = '  +21.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.398'
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.76'
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.45'
$ws.Range("E29").Value = '  +4.28%  '
$ws.Range("D30").Value = '1.912.84'
$ws.Range("E30").Value = '  +3.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.55'
$ws.Range("E31").Value = '  +4.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.117'
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.978'
$ws.Range("E33").Value = '  +5.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08573'
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.690'
$ws.Range("E36").Value = '  +5.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.361'
$ws.Range("E37").Value = '  +3.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06209'
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2152'
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.421'
$ws.Range("E41").Value = '  +2.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.218'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6190'
$ws.Range("E43").Value = '  +4.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.28'
$ws.Range("E44").Value = '  +6.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9977'
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.914'
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6017'
$ws.Range("E47").Value = '  +6.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.33'
$ws.Range("E48").Value = '  +2.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.042'
$ws.Range("E49").Value = '  +4.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07177'
$ws.Range("E50").Value = '  +2.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '76.88'
$ws.Range("E51").Value = '  +2.36%  '
